$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.000006198883056640625
$ws.Range("B3").Value = 0.000005483627319335938
$ws.Range("B4").Value = 0.0004391670227050781
$ws.Range("B5").Value = 0.06713366508483887
$ws.Range("B6").Value = 0.6586453914642334
$ws.Range("B7").Value = 2.751649141311646
$ws.Range("B8").Value = 5.469060659408569
$ws.Range("B9").Value = 17.134361743927
